# Adjust rf to the same period of return (previously annual rf was used in
# all circumstances, which was not correct). Updates GRS, A|a|, A|a|/A|re|,
# and A(a^2)/A(re^2) columns (B, D, E, F) for rows 2-6 on Sheet1.
# Column C (p-value of GRS) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 12.62555539199553
$ws.Range("D2").Value = 0.009911156927268711
$ws.Range("E2").Value = 1.452258593099055
$ws.Range("F2").Value = 2.109055021230046

$ws.Range("B3").Value = 12.65320434403941
$ws.Range("D3").Value = 0.009900073611160558
$ws.Range("E3").Value = 1.450634580768687
$ws.Range("F3").Value = 2.104340686921944

$ws.Range("B4").Value = 13.25310249097062
$ws.Range("D4").Value = 0.009406730847589709
$ws.Range("E4").Value = 1.378346222003244
$ws.Range("F4").Value = 1.899838307710617

$ws.Range("B5").Value = 13.27123803194278
$ws.Range("D5").Value = 0.00984932230049297
$ws.Range("E5").Value = 1.44319811017611
$ws.Range("F5").Value = 2.082820785215896

$ws.Range("B6").Value = 13.26180234133492
$ws.Range("D6").Value = 0.0093988266816417
$ws.Range("E6").Value = 1.37718804309401
$ws.Range("F6").Value = 1.896646906041108

$wb.Save()
